# C5-PowerPoint.pptx edit
#  1) Re-style the "Sources of finance" table on slide 6 to a different
#     built-in PowerPoint table style.
#  2) Re-colour the (single reachable) theme colour scheme from the
#     "Integral" palette to the stock "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
if ($tableShape.HasTable) {
    $tableShape.Table.ApplyStyle("{2C8F72F1-5AA0-47AB-A9EC-D2C626D3494F}")
}

# --- 2. Theme colours -------------------------------------------------
# Helper: "RRGGBB" -> OLE/VBA BGR-packed long used by ThemeColor.RGB
function ConvertTo-OleColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Index order exposed by ThemeColorScheme: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = ConvertTo-OleColor $officeThemeColors[$i - 1]
}
